# Re-sort the alaA mutation table by the "aa position" column (C), ascending,
# treating the values as text (they are stored as text in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dataRange = $ws.Range("A2:D41")
$sortKey = $ws.Range("C2:C41")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($ws.Range("A1:D41"))
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

# Restore view: scroll so row 31 is the top-left visible row, and select F5.
$ws.Range("F5").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
